$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.886.36'
$ws.Range("E2").Value = '  -0.87%  '

$ws.Range("D3").Value = '3.411.13'
$ws.Range("E3").Value = '  -0.51%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '408.15'
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("D6").Value = '128.66'
$ws.Range("E6").Value = '  -1.44%  '

$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +5.92%  '

$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.731'
$ws.Range("E9").Value = '  +5.74%  '

$ws.Range("E10").Value = '  +1.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.62'
$ws.Range("E11").Value = '  +1.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000216'
$ws.Range("E12").Value = '  +45.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.13'
$ws.Range("E13").Value = '  +8.87%  '

$ws.Range("E14").Value = '  -0.26%  '

$ws.Range("D15").Value = '3.961.82'
$ws.Range("E15").Value = '  -0.21%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.22'
$ws.Range("E16").Value = '  +6.93%  '

$ws.Range("D17").Value = '3.399.86'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.46'
$ws.Range("E18").Value = '  +7.68%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.08'
$ws.Range("E19").Value = '  +6.40%  '

$ws.Range("D20").Value = '61.872.31'
$ws.Range("E20").Value = '  -0.72%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '453.18'
$ws.Range("E21").Value = '  +45.41%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '91.44'
$ws.Range("E22").Value = '  +8.25%  '

$ws.Range("E23").Value = '  +1.04%  '

$ws.Range("E24").Value = '  +2.13%  '

$ws.Range("E25").Value = '  +3.81%  '

$ws.Range("B26").Value = 'Filecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.29'
$ws.Range("E26").Value = '  +14.68%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '32.98'
$ws.Range("E27").Value = '  +11.21%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.78'
$ws.Range("E28").Value = '  +0.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.59'
$ws.Range("E29").Value = '  -2.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.69'
$ws.Range("E30").Value = '  -1.67%  '

$ws.Range("D31").Value = '12.02'
$ws.Range("E31").Value = '  +5.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.171'
$ws.Range("E32").Value = '  -0.89%  '

$ws.Range("E33").Value = '  -0.90%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.70'
$ws.Range("E34").Value = '  -4.61%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("E36").Value = '  +3.43%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '54.17'
$ws.Range("E37").Value = '  +4.66%  '

$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  +0.02%  '

$ws.Range("E39").Value = '  +2.16%  '

$ws.Range("E40").Value = '  +7.00%  '

$ws.Range("E41").Value = '  -1.18%  '

$ws.Range("D42").Value = '0.318'
$ws.Range("E42").Value = '  -1.02%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '142.35'
$ws.Range("E43").Value = '  -0.29%  '

$ws.Range("E44").Value = '  +8.71%  '

$ws.Range("E45").Value = '  +14.89%  '

$ws.Range("E46").Value = '  +0.92%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.55'
$ws.Range("E47").Value = '  -1.68%  '

$ws.Range("D48").Value = '22.23'
$ws.Range("E48").Value = '  +5.45%  '

$ws.Range("E49").Value = '  +19.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.14'
$ws.Range("E50").Value = '  +8.24%  '

$ws.Range("D51").Value = '3.764.03'
$ws.Range("E51").Value = '  -0.24%  '
